$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '41.346.23'
$ws.Range('E2').Value = '  -1.09%  '
# Row 3
$ws.Range('D3').Value = '2.189.24'
$ws.Range('E3').Value = '  -1.27%  '
# Row 4
$ws.Range('E4').Value = '  +0.01%  '
# Row 5
$ws.Range('D5').Value = '''253.59'
$ws.Range('E5').Value = '  +2.87%  '
# Row 6
$ws.Range('D6').Value = '''0.624'
$ws.Range('E6').Value = '  -0.68%  '
# Row 7
$ws.Range('D7').Value = '''68.03'
$ws.Range('E7').Value = '  -3.13%  '
# Row 8
$ws.Range('E8').Value = '  +0.07%  '
# Row 9
$ws.Range('D9').Value = '''0.578'
$ws.Range('E9').Value = '  +5.11%  '
# Row 10
$ws.Range('D10').Value = '''37.73'
$ws.Range('E10').Value = '  +2.19%  '
# Row 11
$ws.Range('D11').Value = '''58.27'
$ws.Range('E11').Value = '  +0.36%  '
# Row 12
$ws.Range('D12').Value = '''0.0941'
$ws.Range('E12').Value = '  -1.26%  '
# Row 13
$ws.Range('D13').Value = '''7.04'
$ws.Range('E13').Value = '  +4.69%  '
# Row 14
$ws.Range('D14').Value = '''0.103'
$ws.Range('E14').Value = '  -1.84%  '
# Row 15
$ws.Range('D15').Value = '2.516.02'
$ws.Range('E15').Value = '  -1.24%  '
# Row 16
$ws.Range('D16').Value = '''0.869'
$ws.Range('E16').Value = '  +2.73%  '
# Row 17
$ws.Range('D17').Value = '''14.42'
$ws.Range('E17').Value = '  -2.95%  '
# Row 18
$ws.Range('D18').Value = '2.182.63'
$ws.Range('E18').Value = '  -1.60%  '
# Row 19
$ws.Range('D19').Value = '41.252.09'
$ws.Range('E19').Value = '  -1.15%  '
# Row 20
$ws.Range('D20').Value = '0.0₃0952'
$ws.Range('E20').Value = '  -0.53%  '
# Row 21
$ws.Range('D21').Value = '''6.24'
$ws.Range('E21').Value = '  +2.25%  '
# Row 22
$ws.Range('D22').Value = '''72.19'
$ws.Range('E22').Value = '  -1.67%  '
# Row 23
$ws.Range('D23').Value = '''232.54'
$ws.Range('E23').Value = '  -0.99%  '
# Row 24
$ws.Range('D24').Value = '''2.08'
$ws.Range('E24').Value = '  +0.01%  '
# Row 25
$ws.Range('D25').Value = '''11.95'
$ws.Range('E25').Value = '  +20.39%  '
# Row 26
$ws.Range('D26').Value = '''3.83'
$ws.Range('E26').Value = '  +6.04%  '
# Row 27
$ws.Range('E27').Value = '  +0.04%  '
# Row 28
$ws.Range('D28').Value = '''2.51'
$ws.Range('E28').Value = '  +2.44%  '
# Row 29
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = '''3.73'
$ws.Range('E29').Value = '  -3.61%  '
# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.19'
$ws.Range('E30').Value = '  -1.66%  '
# Row 31
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '''170.08'
$ws.Range('E31').Value = '  -0.10%  '
# Row 32
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''20.59'
$ws.Range('E32').Value = '  +0.49%  '
# Row 33
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '''0.118'
$ws.Range('E33').Value = '  -1.10%  '
# Row 34
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '''0.123'
$ws.Range('E34').Value = '  -1.72%  '
# Row 35
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''5.45'
$ws.Range('E35').Value = '  +6.08%  '
# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.0727'
$ws.Range('E36').Value = '  +1.57%  '
# Row 37
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').Value = '''4.59'
$ws.Range('E37').Value = '  -0.67%  '
# Row 38
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '''25.36'
$ws.Range('E38').Value = '  +7.70%  '
# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''4.03'
$ws.Range('E39').Value = '  +3.52%  '
# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.0299'
$ws.Range('E40').Value = '  +8.92%  '
# Row 41
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '''2.23'
$ws.Range('E41').Value = '  -2.41%  '
# Row 42
$ws.Range('D42').Value = '''12.26'
$ws.Range('E42').Value = '  +18.71%  '
# Row 43
$ws.Range('B43').Value = 'THORChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D43').Value = '''5.73'
$ws.Range('E43').Value = '  -2.62%  '
# Row 44
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = '''64.18'
$ws.Range('E44').Value = '  -1.91%  '
# Row 45
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '''0.202'
$ws.Range('E45').Value = '  +4.57%  '
# Row 46
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = '''4.81'
$ws.Range('E46').Value = '  -2.05%  '
# Row 47
$ws.Range('D47').Value = '''0.102'
# Row 48
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = '''8.62'
$ws.Range('E48').Value = '  -4.19%  '
# Row 49
$ws.Range('B49').Value = 'BinanceUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D49').Value = '''1.00'
$ws.Range('E49').Value = '  +0.26%  '
# Row 50
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = '''1.14'
$ws.Range('E50').Value = '  +3.70%  '
# Row 51
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '''1.17'
$ws.Range('E51').Value = '  -1.04%  '
